$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace test data in rows 3-4 with real patient data ---
$ws.Range("A3").Value = 209315647
$ws.Range("B3").Value = "fadi"
$ws.Range("C3").Value = "badarni"
$ws.Range("D3").Value = 23
$ws.Range("E3").Value = 82
$ws.Range("F3").Value = 184
$ws.Range("G3").Value = 524183083
$ws.Range("H3").Value = "O+"

$ws.Range("A4").Value = 123456789
$ws.Range("B4").Value = "ahmad"
$ws.Range("C4").Value = "sh"
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 80
$ws.Range("F4").Value = 180
$ws.Range("G4").Value = 15241830
$ws.Range("H4").Value = "B"

# --- New "Gender" column (column I) for existing rows ---
$ws.Range("I1").Value = "male"
$ws.Range("I2").Value = "male"
$ws.Range("I3").Value = "male"
$ws.Range("I4").Value = "male"

# --- New row 5: newly diagnosed patient ---
$ws.Range("A5").Value = 123435234
$ws.Range("B5").Value = "abed"
$ws.Range("C5").Value = "ak"
$ws.Range("D5").Value = 23
$ws.Range("E5").Value = 79
$ws.Range("F5").Value = 176
$ws.Range("G5").Value = 524111123
$ws.Range("H5").Value = "A"
$ws.Range("I5").Value = "Male"

[void]$ws.Range("I5").Select()
